# Updated cryptos list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.225.20"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "'1.882.88"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("D4").Value = "'1.009"
$ws.Range("E4").Value = "  +0.54%  "
$ws.Range("D5").Value = "'315.43"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("D7").Value = "'0.5140"
$ws.Range("E7").Value = "  +0.81%  "
$ws.Range("D8").Value = "'0.3901"
$ws.Range("E8").Value = "  +1.40%  "
$ws.Range("D9").Value = "'0.08380"
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("D10").Value = "'1.124"
$ws.Range("E10").Value = "  +1.10%  "
$ws.Range("D11").Value = "'41.66"
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").Value = "'1.880.31"
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").Value = "'7.266"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").Value = "'1.009"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("D17").Value = "'0.00001105"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").Value = "'91.05"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").Value = "'0.06683"
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("D20").Value = "'17.81"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("D21").Value = "'1.008"
$ws.Range("D22").Value = "'6.030"
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("D23").Value = "'28.262.73"
$ws.Range("E23").Value = "  +0.91%  "
$ws.Range("E24").Value = "  +1.00%  "
$ws.Range("D25").Value = "'2.278"
$ws.Range("E25").Value = "  +1.68%  "
$ws.Range("D26").Value = "'2.102.35"
$ws.Range("E26").Value = "  +1.56%  "
$ws.Range("D27").Value = "'160.17"
$ws.Range("E27").Value = "  +1.54%  "
$ws.Range("D28").Value = "'2.471"
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("D30").Value = "'125.52"
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("D31").Value = "'0.1060"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").Value = "'1.042"
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("D33").Value = "'5.874"
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("D34").Value = "'3.616"
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("D35").Value = "'9.465"
$ws.Range("E35").Value = "  +1.10%  "
$ws.Range("E36").Value = "  +1.47%  "
$ws.Range("D37").Value = "'0.06575"
$ws.Range("E37").Value = "  +1.27%  "
$ws.Range("D38").Value = "'0.2214"
$ws.Range("E38").Value = "  +1.78%  "
$ws.Range("D39").Value = "'1.198"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").Value = "'0.6500"
$ws.Range("E40").Value = "  -1.68%  "
$ws.Range("D41").Value = "'1.245"
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("D42").Value = "'5.006"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("E43").Value = "  +0.66%  "
$ws.Range("D44").Value = "'0.6102"
$ws.Range("E44").Value = "  -1.04%  "
$ws.Range("D45").Value = "'13.18"
$ws.Range("E45").Value = "  +1.18%  "
$ws.Range("D46").Value = "'3.695"
$ws.Range("E46").Value = "  +1.22%  "
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").Value = "'2.017"
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("D49").Value = "'1.237"
$ws.Range("E49").Value = "  +2.48%  "
$ws.Range("D50").Value = "'121.09"
$ws.Range("E50").Value = "  +1.18%  "

# Row 51 special-case: coin changed from Cronos to Aave
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'78.24"
$ws.Range("E51").Value = "  -0.59%  "

